$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '288.30'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.16%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.46%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.931'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.43%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07331'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.63%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.215'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '21.49%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.717'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.41%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.97%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.34%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09275'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '19.85%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1696'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.76%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08213'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.06%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03119'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.95%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09945'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.67%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001494'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.32%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005725'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.08%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.535'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.93%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.081'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.02%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.37%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.154'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.96%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-11.95%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04541'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.07%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001210'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.38%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004163'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '3.65%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.98%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003396'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-0.85%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04449'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.06%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007348'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.62%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009532'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-5.41%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1329'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '1.68%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002289'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '13.70%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009084'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.43%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006115'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.15%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.02%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.380'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '5.98%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002099'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.02%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.02%'
